$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: copy formats for new data rows 872:884 from row 871 template ---
$ws.Range("A871:H871").Copy()
$ws.Range("A872:H884").PasteSpecial(-4122)
$ws.Range("I871:J871").Copy()
$ws.Range("I872:J884").PasteSpecial(-4122)

# --- Step 2: copy formats for trailing blank rows 885:901 from a uniformly-styled row template ---
$ws.Range("A719:J719").Copy()
$ws.Range("A885:J901").PasteSpecial(-4122)

# --- Step 3: populate values + formulas for the new data rows ---
$ws.Range("A872").Value = "2023-10-09"
$ws.Range("B872").Value = "18:15"
$ws.Range("C872").Value = "20:00"
$ws.Range("D872").Value = "1h 45m"
$ws.Range("E872").Value = "#python"
$ws.Range("F872").Value = "nwtimetrackingmanager v1.0.0"
$ws.Range("G872").Value = "'True"
$ws.Range("H872").Value = "'False"
$ws.Range("I872").Formula = "=YEAR(A872)"
$ws.Range("J872").Formula = "=MONTH(A872)"
$ws.Range("A873").Value = "2023-10-10"
$ws.Range("B873").Value = "17:45"
$ws.Range("C873").Value = "18:00"
$ws.Range("D873").Value = "0h 15m"
$ws.Range("E873").Value = "#python"
$ws.Range("F873").Value = "nwtimetrackingmanager v1.0.0"
$ws.Range("G873").Value = "'True"
$ws.Range("H873").Value = "'False"
$ws.Range("I873").Formula = "=YEAR(A873)"
$ws.Range("J873").Formula = "=MONTH(A873)"
$ws.Range("A874").Value = "2023-10-11"
$ws.Range("B874").Value = "18:00"
$ws.Range("C874").Value = "19:30"
$ws.Range("D874").Value = "1h 30m"
$ws.Range("E874").Value = "#python"
$ws.Range("F874").Value = "nwtimetrackingmanager v1.0.0"
$ws.Range("G874").Value = "'True"
$ws.Range("H874").Value = "'False"
$ws.Range("I874").Formula = "=YEAR(A874)"
$ws.Range("J874").Formula = "=MONTH(A874)"
$ws.Range("A875").Value = "2023-10-16"
$ws.Range("B875").Value = "08:00"
$ws.Range("C875").Value = "08:30"
$ws.Range("D875").Value = "0h 30m"
$ws.Range("E875").Value = "#python"
$ws.Range("F875").Value = "nwtimetrackingmanager v1.0.0"
$ws.Range("G875").Value = "'True"
$ws.Range("H875").Value = "'False"
$ws.Range("I875").Formula = "=YEAR(A875)"
$ws.Range("J875").Formula = "=MONTH(A875)"
$ws.Range("A876").Value = "2023-10-16"
$ws.Range("B876").Value = "19:00"
$ws.Range("C876").Value = "20:00"
$ws.Range("D876").Value = "1h 00m"
$ws.Range("E876").Value = "#python"
$ws.Range("F876").Value = "nwtimetrackingmanager v1.0.0"
$ws.Range("G876").Value = "'True"
$ws.Range("H876").Value = "'False"
$ws.Range("I876").Formula = "=YEAR(A876)"
$ws.Range("J876").Formula = "=MONTH(A876)"
$ws.Range("A877").Value = "2023-10-17"
$ws.Range("B877").Value = "08:00"
$ws.Range("C877").Value = "08:30"
$ws.Range("D877").Value = "0h 30m"
$ws.Range("E877").Value = "#python"
$ws.Range("F877").Value = "nwtimetrackingmanager v1.0.0"
$ws.Range("G877").Value = "'True"
$ws.Range("H877").Value = "'False"
$ws.Range("I877").Formula = "=YEAR(A877)"
$ws.Range("J877").Formula = "=MONTH(A877)"
$ws.Range("A878").Value = "2023-10-17"
$ws.Range("B878").Value = "17:15"
$ws.Range("C878").Value = "17:45"
$ws.Range("D878").Value = "0h 30m"
$ws.Range("E878").Value = "#python"
$ws.Range("F878").Value = "nwtimetrackingmanager v1.0.0"
$ws.Range("G878").Value = "'True"
$ws.Range("H878").Value = "'False"
$ws.Range("I878").Formula = "=YEAR(A878)"
$ws.Range("J878").Formula = "=MONTH(A878)"
$ws.Range("A879").Value = "2023-10-18"
$ws.Range("B879").Value = "08:00"
$ws.Range("C879").Value = "08:30"
$ws.Range("D879").Value = "0h 30m"
$ws.Range("E879").Value = "#python"
$ws.Range("F879").Value = "nwtimetrackingmanager v1.0.0"
$ws.Range("G879").Value = "'True"
$ws.Range("H879").Value = "'False"
$ws.Range("I879").Formula = "=YEAR(A879)"
$ws.Range("J879").Formula = "=MONTH(A879)"
$ws.Range("A880").Value = "2023-10-18"
$ws.Range("B880").Value = "17:30"
$ws.Range("C880").Value = "17:45"
$ws.Range("D880").Value = "0h 15m"
$ws.Range("E880").Value = "#python"
$ws.Range("F880").Value = "nwtimetrackingmanager v1.0.0"
$ws.Range("G880").Value = "'True"
$ws.Range("H880").Value = "'False"
$ws.Range("I880").Formula = "=YEAR(A880)"
$ws.Range("J880").Formula = "=MONTH(A880)"
$ws.Range("A881").Value = "2023-10-19"
$ws.Range("B881").Value = "08:00"
$ws.Range("C881").Value = "08:30"
$ws.Range("D881").Value = "0h 30m"
$ws.Range("E881").Value = "#python"
$ws.Range("F881").Value = "nwtimetrackingmanager v1.0.0"
$ws.Range("G881").Value = "'True"
$ws.Range("H881").Value = "'False"
$ws.Range("I881").Formula = "=YEAR(A881)"
$ws.Range("J881").Formula = "=MONTH(A881)"
$ws.Range("A882").Value = "2023-10-19"
$ws.Range("B882").Value = "17:30"
$ws.Range("C882").Value = "17:45"
$ws.Range("D882").Value = "0h 15m"
$ws.Range("E882").Value = "#python"
$ws.Range("F882").Value = "nwtimetrackingmanager v1.0.0"
$ws.Range("G882").Value = "'True"
$ws.Range("H882").Value = "'False"
$ws.Range("I882").Formula = "=YEAR(A882)"
$ws.Range("J882").Formula = "=MONTH(A882)"
$ws.Range("A883").Value = "2023-10-20"
$ws.Range("B883").Value = "08:00"
$ws.Range("C883").Value = "08:30"
$ws.Range("D883").Value = "0h 30m"
$ws.Range("E883").Value = "#python"
$ws.Range("F883").Value = "nwtimetrackingmanager v1.0.0"
$ws.Range("G883").Value = "'True"
$ws.Range("H883").Value = "'False"
$ws.Range("I883").Formula = "=YEAR(A883)"
$ws.Range("J883").Formula = "=MONTH(A883)"
$ws.Range("A884").Value = "2023-10-20"
$ws.Range("B884").Value = "20:00"
$ws.Range("C884").Value = "00:00"
$ws.Range("D884").Value = "4h 00m"
$ws.Range("E884").Value = "#python"
$ws.Range("F884").Value = "nwtimetrackingmanager v1.0.0"
$ws.Range("G884").Value = "'True"
$ws.Range("H884").Value = "'False"
$ws.Range("I884").Formula = "=YEAR(A884)"
$ws.Range("J884").Formula = "=MONTH(A884)"

# --- Step 4: update dimension / pane / selection to reflect new extent ---
$ws.Range("A1").Select()
$ws.Application.ActiveWindow.ScrollRow = 864
$ws.Range("E887").Select()
